# Auto-generated COM-interop script applying the scheduled market-data refresh
# described by the commit "chore: update Sheets via scheduled runner".
# For each leve row touched by the refresh we rewrite the price/profit columns
# (H..N) with the freshly pulled Universalis averages, matching the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 71431760
$ws.Range("I64").Value = 250002030
$ws.Range("J64").Value = 3656
$ws.Range("K64").Value = 250002030
$ws.Range("L64").Value = 3656
$ws.Range("M64").Value = -250001782
$ws.Range("N64").Value = -4152
# Row 67
$ws.Range("H67").Value = 71431760
$ws.Range("I67").Value = 250002030
$ws.Range("J67").Value = 3656
$ws.Range("K67").Value = 250002030
$ws.Range("L67").Value = 3656
$ws.Range("M67").Value = -250001172
$ws.Range("N67").Value = -5372
# Row 70
$ws.Range("H70").Value = 4412.625
$ws.Range("I70").Value = 4800.1665
$ws.Range("J70").Value = 3250
$ws.Range("K70").Value = 14400.4995
$ws.Range("L70").Value = 9750
$ws.Range("M70").Value = -14130.4995
$ws.Range("N70").Value = -10290
# Row 73
$ws.Range("H73").Value = 4412.625
$ws.Range("I73").Value = 4800.1665
$ws.Range("J73").Value = 3250
$ws.Range("K73").Value = 14400.4995
$ws.Range("L73").Value = 9750
$ws.Range("M73").Value = -13464.4995
$ws.Range("N73").Value = -11622
# Row 74
$ws.Range("H74").Value = 3383.923
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 3499.2173
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 3499.2173
$ws.Range("M74").Value = -1564
$ws.Range("N74").Value = -5371.2173
# Row 77
$ws.Range("H77").Value = 3383.923
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 3499.2173
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 17496.0865
$ws.Range("M77").Value = -7820
$ws.Range("N77").Value = -26856.0865
# Row 86
$ws.Range("H86").Value = 1960
$ws.Range("I86").Value = 1825
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1825
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -702
$ws.Range("N86").Value = -4746
# Row 89
$ws.Range("H89").Value = 1960
$ws.Range("I89").Value = 1825
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 9125
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -3509
$ws.Range("N89").Value = -23732
# Row 132
$ws.Range("H132").Value = 2476.2646
$ws.Range("I132").Value = 1840.52
$ws.Range("J132").Value = 4242.222
$ws.Range("K132").Value = 5521.559999999999
$ws.Range("L132").Value = 12726.666
$ws.Range("M132").Value = -2991.559999999999
$ws.Range("N132").Value = -17786.666
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 137
$ws.Range("H137").Value = 2941.6987
$ws.Range("I137").Value = 1101.7941
$ws.Range("J137").Value = 4545.718
$ws.Range("K137").Value = 3305.3823
$ws.Range("L137").Value = 13637.154
$ws.Range("M137").Value = -755.3823000000002
$ws.Range("N137").Value = -18737.154
# Row 138
$ws.Range("H138").Value = 1593.2
$ws.Range("I138").Value = 901.2121
$ws.Range("J138").Value = 2936.4707
$ws.Range("K138").Value = 2703.6363
$ws.Range("L138").Value = 8809.4121
$ws.Range("M138").Value = 2436.3637
$ws.Range("N138").Value = -19089.4121

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 851.5
$ws.Range("I3").Value = 800
$ws.Range("K3").Value = 800
$ws.Range("M3").Value = -685
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 17
$ws.Range("H17").Value = 1500
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1846
# Row 122
$ws.Range("H122").Value = 2078.5454
$ws.Range("I122").Value = 2281.25
$ws.Range("J122").Value = 1538
$ws.Range("K122").Value = 6843.75
$ws.Range("L122").Value = 4614
$ws.Range("M122").Value = -4393.75
$ws.Range("N122").Value = -9514

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 336.33334
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 502.5
$ws.Range("K8").Value = 4
$ws.Range("L8").Value = 502.5
$ws.Range("M8").Value = 136
$ws.Range("N8").Value = -782.5
# Row 10
$ws.Range("H10").Value = 1003
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1003
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1003
$ws.Range("N10").Value = -1283
$ws.Range("M10").ClearContents()
# Row 18
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -11058

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 12064501
$ws.Range("I31").Value = 37037972
$ws.Range("J31").Value = 23720
$ws.Range("K31").Value = 37037972
$ws.Range("L31").Value = 23720
$ws.Range("M31").Value = -37037677
$ws.Range("N31").Value = -24310
# Row 34
$ws.Range("H34").Value = 12064501
$ws.Range("I34").Value = 37037972
$ws.Range("J34").Value = 23720
$ws.Range("K34").Value = 37037972
$ws.Range("L34").Value = 23720
$ws.Range("M34").Value = -37037770
$ws.Range("N34").Value = -24124

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 2181.7646
$ws.Range("J55").Value = 2576.923
$ws.Range("L55").Value = 7730.768999999999
$ws.Range("N55").Value = -8084.768999999999
# Row 122
$ws.Range("H122").Value = 862.89746
$ws.Range("I122").Value = 392.0435
$ws.Range("J122").Value = 1539.75
$ws.Range("K122").Value = 3528.3915
$ws.Range("L122").Value = 13857.75
$ws.Range("M122").Value = -1078.3915
$ws.Range("N122").Value = -18757.75
# Row 131
$ws.Range("H131").Value = 955.90625
$ws.Range("I131").Value = 290
$ws.Range("J131").Value = 1024.7931
$ws.Range("K131").Value = 870
$ws.Range("L131").Value = 3074.379300000001
$ws.Range("M131").Value = 4170
$ws.Range("N131").Value = -13154.3793
# Row 137
$ws.Range("H137").Value = 4255717.5
$ws.Range("I137").Value = 62730.555
$ws.Range("J137").Value = 16834678
$ws.Range("K137").Value = 188191.665
$ws.Range("L137").Value = 50504034
$ws.Range("M137").Value = -183091.665
$ws.Range("N137").Value = -50514234
# Row 139
$ws.Range("H139").Value = 409667.28
$ws.Range("I139").Value = 1001042.75
$ws.Range("J139").Value = 3096.6875
$ws.Range("K139").Value = 3003128.25
$ws.Range("L139").Value = 9290.0625
$ws.Range("M139").Value = -2997988.25
$ws.Range("N139").Value = -19570.0625

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2312.5557
$ws.Range("I80").Value = 1976.625
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 1976.625
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -978.625
$ws.Range("N80").Value = -6996
# Row 83
$ws.Range("H83").Value = 2312.5557
$ws.Range("I83").Value = 1976.625
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 9883.125
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -4891.125
$ws.Range("N83").Value = -34984
# Row 122
$ws.Range("H122").Value = 2426.6572
$ws.Range("I122").Value = 1920.36
$ws.Range("J122").Value = 3692.4
$ws.Range("K122").Value = 5761.08
$ws.Range("L122").Value = 11077.2
$ws.Range("M122").Value = -3311.08
$ws.Range("N122").Value = -15977.2
# Row 132
$ws.Range("H132").Value = 5533.7715
$ws.Range("I132").Value = 6019.615
$ws.Range("J132").Value = 4130.222
$ws.Range("K132").Value = 18058.845
$ws.Range("L132").Value = 12390.666
$ws.Range("M132").Value = -15528.845
$ws.Range("N132").Value = -17450.666

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1287.1428
$ws.Range("I122").Value = 1192
$ws.Range("J122").Value = 1525
$ws.Range("K122").Value = 3576
$ws.Range("L122").Value = 4575
$ws.Range("M122").Value = -1126
$ws.Range("N122").Value = -9475
# Row 136
$ws.Range("H136").Value = 4203.2905
$ws.Range("I136").Value = 4307.2964
$ws.Range("J136").Value = 3501.25
$ws.Range("K136").Value = 12921.8892
$ws.Range("L136").Value = 10503.75
$ws.Range("M136").Value = -10371.8892
$ws.Range("N136").Value = -15603.75
# Row 137
$ws.Range("H137").Value = 46528.832
$ws.Range("J137").Value = 46528.832
$ws.Range("L137").Value = 46528.832
$ws.Range("N137").Value = -56728.832
# Row 138
$ws.Range("H138").Value = 48907.25
$ws.Range("J138").Value = 48907.25
$ws.Range("L138").Value = 48907.25
$ws.Range("N138").Value = -59187.25
# Row 139
$ws.Range("H139").Value = 37707.5
$ws.Range("J139").Value = 37707.5
$ws.Range("L139").Value = 37707.5
$ws.Range("N139").Value = -47987.5
